# DatDF now has two level index and passes unittests
# Rebuild the sheet content: insert a new "datnum"/"datname" pair of
# columns, add a second data row ("base" level) and keep the original
# data row (now row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember a cell that already carries the bold/bordered header style
# (style index 1 in the original workbook) so we can stamp it onto the
# cells that need it before we overwrite their contents.
$ws.Range("B1").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)
$ws.Range("A2:B2").PasteSpecial(-4122)
$ws.Range("A3:B3").PasteSpecial(-4122)

# Now clear any left-over values (formats were already applied above and
# are unaffected by ClearContents).
$ws.Cells.ClearContents()

# ---- Row 1 : header ----
$ws.Range("A1").Value = "datnum"
$ws.Range("B1").Value = "datname"
$ws.Range("C1").Value = "time"
$ws.Range("D1").Value = "picklepath"
$ws.Range("E1").Value = "x_label"
$ws.Range("F1").Value = "y_label"
$ws.Range("G1").Value = "dim"
$ws.Range("H1").Value = "time_elapsed"

# ---- Row 2 : new "base" level row ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "base"
$ws.Range("C2").Value = 1577779311.487608
$ws.Range("D2").Value = "pathtopickle"
$ws.Range("E2").Value = "xlabel"
$ws.Range("F2").Value = "ylabel"

# ---- Row 3 : original data row (shifted down one row) ----
$ws.Range("A3").Value = 2700
$ws.Range("B3").Value = "base"
$ws.Range("E3").Value = "FD_SDP/1000mV"
$ws.Range("F3").Value = "Repeats (mV)"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 13.421

$null = $ws.Range("A1").Select()
